# Loan RBI, Variable Instalments
#
# Inserts a new (blank) column before the existing "In Advance" column on
# the "Repayment schedule" sheet, shifting the trailing columns one to the
# right, and moves the active sheet/selection from "Transactions" to
# "Repayment schedule".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column N (pushes former N/O/P -> O/P/Q).
$ws.Columns("N:N").Insert() | Out-Null

# Give the newly inserted column the same width as column M ("Due").
$ws.Columns("N:N").ColumnWidth = 9.8

# Make "Repayment schedule" the active sheet/tab (was "Transactions"),
# matching the new selection left on it.
$ws.Activate() | Out-Null
$ws.Range("S9").Select() | Out-Null
